$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$BValues = @(18.08625015614325, 17.64250469793264, 17.36624223006096, 17.25286273582108, 17.23399251259574, 17.36471618393429, 17.93412649843798, 19.01415583286884, 19.77736692925943, 20.11649095246747, 20.24363613226981, 20.21631150781764, 20.12697741001229, 20.07208859522178, 19.75503213112902, 19.55837856966187, 19.44451578194948, 19.40583812493968, 19.57939145150646, 20.15325243070719, 20.52083321655464, 20.32536723465649, 19.56989402324655, 18.72675575454442)
$CValues = @(9.738717452883067, 9.515025990553843, 9.374347221146962, 9.316240138901613, 9.306546200830699, 9.373566644749017, 9.662314023781253, 10.19972357362072, 10.57403431950758, 10.73930793975763, 10.80113195028808, 10.78785154747384, 10.7444097297568, 10.71769999609069, 10.56312895429614, 10.46699229366441, 10.41122849539168, 10.39226870382453, 10.4772750136375, 10.75719063129969, 10.93567578680527, 10.84083566395279, 10.47262773003469, 10.05774446811624)
$DValues = @(7.997305569775516, 7.999093464023706, 8.000850131173173, 8.001732233848667, 8.001888763253156, 8.000861353655701, 7.997785624862069, 7.996957290825876, 7.999487087164724, 8.001311853207463, 8.002099062000571, 8.001925255031688, 8.001374693286428, 8.001049966327539, 7.99938133642409, 7.998529807932285, 7.998103510236877, 7.99797009675335, 7.998613889783303, 8.001533801250821, 8.004002570610272, 8.002633900706785, 7.998575679298731, 7.996627718330883)
$EValues = @(12.98521812365697, 13.02009814429155, 13.04269451355871, 13.05220020832657, 13.05379661286617, 13.04282150509225, 12.99700034779667, 12.91647073980015, 12.86293971281614, 12.83979961382768, 12.83121044968222, 12.83305257517595, 12.83908950520639, 12.84280987150908, 12.86447628869599, 12.87807770121614, 12.88601494205498, 12.88872196858375, 12.87661800720425, 12.83731160972676, 12.81263350792646, 12.82571241176064, 12.87727756847588, 12.93726306837042)
$FValues = @(36.90534111957079, 37.0406313354654, 37.13270598905642, 37.17248566416525, 37.17922727479143, 37.13323333531011, 36.95011639964445, 36.66274751145453, 36.49569197463769, 36.42933460500373, 36.40559796223476, 36.41064810084859, 36.42735385899771, 36.43776798924566, 36.50022299721449, 36.54100951890916, 36.56537561050349, 36.57378111055954, 36.53657383547846, 36.42240916565146, 36.35590881353632, 36.39065718581891, 36.53857635077513, 36.7327728639119)
$HValues = @(7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261, 7.344005520526261)
$JValues = @(10.22553441222674, 10.24710118614279, 10.26102417996165, 10.26686966170441, 10.26785068826062, 10.26110231806797, 10.23282963271189, 10.1827654314856, 10.14922801773629, 10.13466829144124, 10.12925452808322, 10.130416052837, 10.13422090270672, 10.13656445130471, 10.15019350632211, 10.15873256415344, 10.16370960132937, 10.16540602172612, 10.157816781338, 10.13310062436597, 10.1175280296252, 10.12578642785771, 10.15823059545558, 10.19573686370757)
$LValues = @(10.92129140794003, 10.91766058614524, 10.91645423496703, 10.91622105425582, 10.91619797577651, 10.9164500423757, 10.91982790202492, 10.93451580147236, 10.95014996457139, 10.95829715904485, 10.96152957711779, 10.96082689513907, 10.95856015026626, 10.95719083158072, 10.94963822069058, 10.94526894542937, 10.94285338564991, 10.94205231739906, 10.94572397943472, 10.95922196487213, 10.96890118889101, 10.9636572882052, 10.94551795818409, 10.92968630897107)
$MValues = @(16.81086850035017, 16.70607724895481, 16.64321381675943, 16.61798642928898, 16.61382151245772, 16.64287198890557, 16.77444057034852, 17.04336286279494, 17.24647410129748, 17.33983062985395, 17.37529927681542, 17.3676556291827, 17.34274644454404, 17.32750338783001, 17.24039051792208, 17.18718026079912, 17.156666806557, 17.14635187447148, 17.19283525060277, 17.35005988471353, 17.45348459429822, 17.39823077528665, 17.1902783856304, 16.96955775103271)
$OValues = @(28.45777371119365, 28.57547145604962, 28.65438626626417, 28.68821249252116, 28.69392992849923, 28.65483571090583, 28.49697380729973, 28.24032636897388, 28.08425919176382, 28.02036374992901, 27.99719315660426, 28.00213769526436, 28.01843692046532, 28.02855430434787, 28.08857814715416, 28.12722241004089, 28.15011763520762, 28.15798417662051, 28.12303948903077, 28.01362158553788, 27.94808833788178, 27.98251634024155, 28.12492847558199, 28.30406686319026)

$rows = @(2..25)
$cols = @("B","C","D","E","F","H","J","L","M","O")
$allValues = @{
    "B" = $BValues
    "C" = $CValues
    "D" = $DValues
    "E" = $EValues
    "F" = $FValues
    "H" = $HValues
    "J" = $JValues
    "L" = $LValues
    "M" = $MValues
    "O" = $OValues
}

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $allValues[$c][$i]
    }
}

Write-Host "Updated loading percent values for 380 kV case"
